# Auto-generated edit script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.394.71'
$ws.Range("E2").Value = '  +8.65%  '

$ws.Range("D3").Value = '1.602.08'
$ws.Range("E3").Value = '  +8.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9919'
$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.08'
$ws.Range("E6").Value = '  +8.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3690'
$ws.Range("E7").Value = '  +0.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3407'
$ws.Range("E8").Value = '  +9.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.36'
$ws.Range("E9").Value = '  +5.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("E10").Value = '  +6.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07067'
$ws.Range("E11").Value = '  +5.53%  '

$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.76'
$ws.Range("E13").Value = '  +8.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.943'
$ws.Range("E14").Value = '  +7.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.645'
$ws.Range("E15").Value = '  +6.65%  '

$ws.Range("E16").Value = '  +5.56%  '

$ws.Range("D17").Value = '1.603.35'
$ws.Range("E17").Value = '  +8.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9918'
$ws.Range("E18").Value = '  +2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06838'
$ws.Range("E19").Value = '  +14.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.99'

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.14'
$ws.Range("E21").Value = '  +10.67%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.038'
$ws.Range("E22").Value = '  +9.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.84'
$ws.Range("E23").Value = '  +6.66%  '

$ws.Range("D24").Value = '22.439.80'
$ws.Range("E24").Value = '  +8.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.405'
$ws.Range("E25").Value = '  +5.67%  '

$ws.Range("E26").Value = '  +18.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.01'
$ws.Range("E27").Value = '  +5.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.58'
$ws.Range("E28").Value = '  +12.74%  '

$ws.Range("D29").Value = '1.782.49'
$ws.Range("E29").Value = '  +8.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.05'
$ws.Range("E30").Value = '  +5.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.175'
$ws.Range("E31").Value = '  +5.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.145'
$ws.Range("E32").Value = '  +21.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9526'
$ws.Range("E33").Value = '  +14.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08298'
$ws.Range("E34").Value = '  +3.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.638'
$ws.Range("E35").Value = '  +6.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.293'
$ws.Range("E36").Value = '  +10.58%  '

$ws.Range("E37").Value = '  +14.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.273'
$ws.Range("E38").Value = '  +4.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.618'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06121'
$ws.Range("E40").Value = '  +5.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02221'
$ws.Range("E41").Value = '  +8.18%  '

$ws.Range("E42").Value = '  +7.42%  '

$ws.Range("E43").Value = '  +2.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5916'
$ws.Range("E44").Value = '  +11.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.825'
$ws.Range("E45").Value = '  +7.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.17'
$ws.Range("E46").Value = '  +7.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5697'
$ws.Range("E47").Value = '  +8.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.54'
$ws.Range("E48").Value = '  +6.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.985'
$ws.Range("E49").Value = '  +8.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06809'
$ws.Range("E50").Value = '  +4.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.95'
$ws.Range("E51").Value = '  +8.52%  '
